# Auto-generated Excel COM-interop script
# Applies the 2025-11-27 violent-crime data update across the
# 'Citywide Totals', 'By Neighborhood', and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 6052
$ws.Range('L3').Value = 6594
$ws.Range('E4').Value = 2065
$ws.Range('J4').Value = 1881
$ws.Range('L4').Value = 1617
$ws.Range('L5').Value = 391
$ws.Range('L6').Value = 5406
$ws.Range('E7').Value = 26070
$ws.Range('J7').Value = 29357
$ws.Range('L7').Value = 20060

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L5').Value = 71
$ws.Range('L6').Value = 156
$ws.Range('L7').Value = 646
$ws.Range('L8').Value = 1330
$ws.Range('L15').Value = 162
$ws.Range('L16').Value = 44
$ws.Range('L18').Value = 135
$ws.Range('L19').Value = 541
$ws.Range('L22').Value = 64
$ws.Range('L25').Value = 120
$ws.Range('L29').Value = 1125
$ws.Range('L31').Value = 197
$ws.Range('L33').Value = 901
$ws.Range('L36').Value = 254
$ws.Range('L37').Value = 767
$ws.Range('L40').Value = 53
$ws.Range('L44').Value = 137
$ws.Range('L47').Value = 141
$ws.Range('L48').Value = 264
$ws.Range('L50').Value = 98
$ws.Range('L51').Value = 252
$ws.Range('L53').Value = 221
$ws.Range('L54').Value = 436
$ws.Range('J63').Value = 233
$ws.Range('L65').Value = 393
$ws.Range('L67').Value = 693
$ws.Range('L68').Value = 62
$ws.Range('L75').Value = 73
$ws.Range('L78').Value = 259
$ws.Range('L83').Value = 437
$ws.Range('L85').Value = 995
$ws.Range('L88').Value = 212
$ws.Range('L89').Value = 278
$ws.Range('L90').Value = 208
$ws.Range('L91').Value = 271
$ws.Range('L92').Value = 61
$ws.Range('L94').Value = 249
$ws.Range('L95').Value = 287
$ws.Range('E96').Value = 280
$ws.Range('L96').Value = 223
$ws.Range('L97').Value = 162
$ws.Range('L99').Value = 348
$ws.Range('E101').Value = 26070
$ws.Range('J101').Value = 29357
$ws.Range('L101').Value = 20060

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L3').Value = 63
$ws.Range('E4').Value = 21
$ws.Range('E7').Value = 280
$ws.Range('L7').Value = 223

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 222
$ws.Range('L7').Value = 646

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('L2').Value = 71
$ws.Range('L6').Value = 79
$ws.Range('L7').Value = 278

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 298
$ws.Range('L3').Value = 415
$ws.Range('L6').Value = 206
$ws.Range('L7').Value = 995

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('L6').Value = 75
$ws.Range('L7').Value = 221

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 398
$ws.Range('L3').Value = 471
$ws.Range('L7').Value = 1330

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L3').Value = 177
$ws.Range('L6').Value = 97
$ws.Range('L7').Value = 437

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 243
$ws.Range('L3').Value = 316
$ws.Range('L4').Value = 60
$ws.Range('L7').Value = 901

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L6').Value = 68
$ws.Range('L7').Value = 287

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 232
$ws.Range('L3').Value = 271
$ws.Range('L5').Value = 22
$ws.Range('L7').Value = 767

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L6').Value = 94
$ws.Range('L7').Value = 393

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L2').Value = 102
$ws.Range('L7').Value = 348

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('L3').Value = 52
$ws.Range('L7').Value = 197

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L3').Value = 271
$ws.Range('L6').Value = 161
$ws.Range('L7').Value = 693

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L6').Value = 210
$ws.Range('L7').Value = 436

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 337
$ws.Range('L7').Value = 1125

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L2').Value = 40
$ws.Range('L6').Value = 105
$ws.Range('L7').Value = 264

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 194
$ws.Range('L7').Value = 541

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('L6').Value = 33
$ws.Range('L7').Value = 137

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('L3').Value = 46
$ws.Range('L7').Value = 156

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L3').Value = 223
$ws.Range('L6').Value = 178

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L6').Value = 73
$ws.Range('L7').Value = 259

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L3').Value = 124
$ws.Range('L6').Value = 35
$ws.Range('L7').Value = 271

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('L3').Value = 47
$ws.Range('L7').Value = 135

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L6').Value = 62
$ws.Range('L7').Value = 254

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L4').Value = 32
$ws.Range('L7').Value = 249

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('L3').Value = 57
$ws.Range('L7').Value = 120

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('L2').Value = 53
$ws.Range('L7').Value = 141

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('L3').Value = 53
$ws.Range('L7').Value = 162

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('L6').Value = 26
$ws.Range('L7').Value = 98

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('L2').Value = 39
$ws.Range('L7').Value = 162

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('L6').Value = 22
$ws.Range('L7').Value = 61

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('L3').Value = 72
$ws.Range('L7').Value = 212

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('L6').Value = 31
$ws.Range('L7').Value = 71

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('L3').Value = 25
$ws.Range('L7').Value = 73

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('L3').Value = 60
$ws.Range('L7').Value = 208

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L2').Value = 79
$ws.Range('L6').Value = 53
$ws.Range('L7').Value = 252

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('L6').Value = 17
$ws.Range('L7').Value = 62

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('L6').Value = 11
$ws.Range('L7').Value = 64

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range('L3').Value = 25
$ws.Range('L7').Value = 53

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('L3').Value = 6
$ws.Range('L6').Value = 28
$ws.Range('L7').Value = 44

